# Update B2:B285 to 10 and B286:B347 to 20 on the "Dataset1" sheet,
# and update the sheet's view (topLeftCell / selected cell) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dataset1")

$ws.Range("B2:B285").Value = 10
$ws.Range("B286:B347").Value = 20

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 273
$ws.Range("B286").Select()
